$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "66.722.50"
Set-TextValue $ws "D3" "3.091.00"
Set-TextValue $ws "E3" "  +5.53%  "
Set-TextValue $ws "E4" "  +0.01%  "
Set-TextValue $ws "D5" "579.94"
Set-TextValue $ws "E5" "  +2.24%  "
Set-TextValue $ws "D6" "167.86"
Set-TextValue $ws "E6" "  +6.23%  "
Set-TextValue $ws "E7" "  -0.01%  "
Set-TextValue $ws "D8" "3.087.43"
Set-TextValue $ws "E8" "  +5.57%  "
Set-TextValue $ws "E9" "  +1.81%  "
Set-TextValue $ws "D10" "6.60"
Set-TextValue $ws "E10" "  -1.62%  "
Set-TextValue $ws "E11" "  +3.60%  "
Set-TextValue $ws "D12" "0.481"
Set-TextValue $ws "E12" "  +4.53%  "
Set-TextValue $ws "E13" "  +3.06%  "
Set-TextValue $ws "D14" "36.42"
Set-TextValue $ws "E14" "  +6.45%  "
Set-TextValue $ws "E15" "  -0.47%  "
Set-TextValue $ws "D16" "3.604.40"
Set-TextValue $ws "E16" "  +5.41%  "
Set-TextValue $ws "D17" "66.750.99"
Set-TextValue $ws "E17" "  +2.31%  "
Set-TextValue $ws "D18" "7.19"
Set-TextValue $ws "E18" "  +3.14%  "
Set-TextValue $ws "D19" "3.093.30"
Set-TextValue $ws "E19" "  +5.55%  "
Set-TextValue $ws "D20" "16.21"
Set-TextValue $ws "E20" "  +3.92%  "
Set-TextValue $ws "D21" "466.80"
Set-TextValue $ws "E21" "  +5.41%  "
Set-TextValue $ws "E22" "  +3.61%  "
Set-TextValue $ws "D23" "7.50"
Set-TextValue $ws "E23" "  +3.51%  "
Set-TextValue $ws "D24" "83.95"
Set-TextValue $ws "E24" "  +2.28%  "
Set-TextValue $ws "E25" "  +6.55%  "
Set-TextValue $ws "D26" "13.05"
Set-TextValue $ws "E26" "  +8.20%  "
Set-TextValue $ws "D27" "10.10"
Set-TextValue $ws "E28" "  -0.01%  "
Set-TextValue $ws "D29" "8.02"
Set-TextValue $ws "E29" "  -0.21%  "
Set-TextValue $ws "E30" "  +2.07%  "
Set-TextValue $ws "E31" "  +4.08%  "
Set-TextValue $ws "E32" "  +1.50%  "
Set-TextValue $ws "D33" "28.23"
Set-TextValue $ws "E33" "  +4.37%  "
Set-TextValue $ws "E34" "  +3.53%  "
Set-TextValue $ws "D35" "1.00"
Set-TextValue $ws "E35" "  +0.09%  "
Set-TextValue $ws "E36" "  +3.66%  "
Set-TextValue $ws "D37" "5.89"
Set-TextValue $ws "E37" "  +2.96%  "
Set-TextValue $ws "D38" "47.21"
Set-TextValue $ws "E38" "  +6.01%  "
Set-TextValue $ws "E39" "  +6.59%  "
Set-TextValue $ws "E40" "  +6.66%  "
Set-TextValue $ws "D41" "50.30"
Set-TextValue $ws "E41" "  +1.15%  "
Set-TextValue $ws "E42" "  +1.78%  "
Set-TextValue $ws "D43" "8.68"
Set-TextValue $ws "E43" "  +2.50%  "
Set-TextValue $ws "D44" "2.81"
Set-TextValue $ws "E44" "  -0.48%  "
Set-TextValue $ws "E45" "  +3.11%  "
Set-TextValue $ws "D46" "383.14"
Set-TextValue $ws "E46" "  +0.30%  "
Set-TextValue $ws "D47" "2.784.50"
Set-TextValue $ws "E47" "  +3.18%  "
Set-TextValue $ws "D48" "134.97"
Set-TextValue $ws "E48" "  +0.91%  "
Set-TextValue $ws "D50" "24.91"
Set-TextValue $ws "E50" "  +6.92%  "
Set-TextValue $ws "D51" "2.22"
Set-TextValue $ws "E51" "  +1.59%  "
